# Update countries & provincias Spain
# - Refreshes COVID case/death counters for a set of countries.
# - Five country pairs swapped leaderboard position (their row in the
#   sheet stays fixed, but which country name/data occupies that row
#   changes) because the updated "Casos totales" (col B) reordered the
#   descending ranking: Japon/Polonia, Zambia/Libia, Mozambique/Eslovaquia,
#   Burkina Faso/Letonia, Islas Malvinas/Montserrat.
# - The "Datos actualizados" footer timestamp moves from 10:00 to 11:17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 11:17"

# Row 25 - Filipinas
$ws.Range("B25").Value = 197164
$ws.Range("C25").Value = 2965
$ws.Range("D25").Value = 132396
$ws.Range("E25").Value = 61730
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 3038

# Row 26 - Indonesia
$ws.Range("B26").Value = 157859
$ws.Range("C26").Value = 2447
$ws.Range("D26").Value = 112867
$ws.Range("E26").Value = 38134
$ws.Range("G26").Value = 99
$ws.Range("H26").Value = 6858

# Row 32 - Israel
$ws.Range("B32").Value = 105252
$ws.Range("C32").Value = 780
$ws.Range("D32").Value = 83028
$ws.Range("E32").Value = 21368
$ws.Range("G32").Value = 9
$ws.Range("H32").Value = 856

# Row 47 - was Japon, now Polonia (Polonia overtakes Japon)
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 63073
$ws.Range("C47").Value = 763
$ws.Range("D47").Value = 42784
$ws.Range("E47").Value = 18312
$ws.Range("G47").Value = 17
$ws.Range("H47").Value = 1977

# Row 48 - was Polonia, now Japon
$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 62507
$ws.Range("D48").Value = 49340
$ws.Range("E48").Value = 11986
$ws.Range("H48").Value = 1181

# Row 71 - Austria
$ws.Range("B71").Value = 25706
$ws.Range("C71").Value = 211
$ws.Range("D71").Value = 21888
$ws.Range("E71").Value = 3085

# Row 73 - El Salvador
$ws.Range("B73").Value = 24986
$ws.Range("C73").Value = 175
$ws.Range("D73").Value = 12717
$ws.Range("E73").Value = 11591

# Row 88 - was Zambia, now Libia (Libia overtakes Zambia)
$ws.Range("A88").Value = "Libia"
$ws.Range("B88").Value = 11281
$ws.Range("C88").Value = 272
$ws.Range("D88").Value = 1112
$ws.Range("E88").Value = 9966
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 203

# Row 89 - was Libia, now Zambia
$ws.Range("A89").Value = "Zambia"
$ws.Range("B89").Value = 11148
$ws.Range("D89").Value = 10208
$ws.Range("E89").Value = 660
$ws.Range("H89").Value = 280

# Row 92 - Malasia
$ws.Range("B92").Value = 9285
$ws.Range("C92").Value = 11
$ws.Range("D92").Value = 8971
$ws.Range("E92").Value = 189

# Row 111 - Hong Kong
$ws.Range("B111").Value = 4711
$ws.Range("C111").Value = 19
$ws.Range("D111").Value = 4108
$ws.Range("E111").Value = 526

# Row 120 - was Mozambique, now Eslovaquia (Eslovaquia overtakes Mozambique)
$ws.Range("A120").Value = "Eslovaquia"
$ws.Range("B120").Value = 3452
$ws.Range("C120").Value = 28
$ws.Range("D120").Value = 2167
$ws.Range("E120").Value = 1252
$ws.Range("H120").Value = 33

# Row 121 - was Eslovaquia, now Mozambique
$ws.Range("A121").Value = "Mozambique"
$ws.Range("B121").Value = 3440
$ws.Range("D121").Value = 1661
$ws.Range("E121").Value = 1758
$ws.Range("H121").Value = 21

# Row 129 - Lituania
$ws.Range("B129").Value = 2694
$ws.Range("C129").Value = 21
$ws.Range("D129").Value = 1785
$ws.Range("E129").Value = 824

# Row 130 - Eslovenia
$ws.Range("B130").Value = 2686
$ws.Range("C130").Value = 21
$ws.Range("D130").Value = 2139
$ws.Range("E130").Value = 414

# Row 152 - was Burkina Faso, now Letonia (Letonia overtakes Burkina Faso)
$ws.Range("A152").Value = "Letonia"
$ws.Range("B152").Value = 1342
$ws.Range("C152").Value = 5
$ws.Range("D152").Value = 1135
$ws.Range("E152").Value = 174
$ws.Range("H152").Value = 33

# Row 153 - was Letonia, now Burkina Faso
$ws.Range("A153").Value = "Burkina Faso"
$ws.Range("B153").Value = 1338
$ws.Range("D153").Value = 1050
$ws.Range("E153").Value = 233
$ws.Range("H153").Value = 55

# Row 191 - Brunei
$ws.Range("B191").Value = 144
$ws.Range("C191").Value = 1
$ws.Range("E191").Value = 2

# Row 214 - was Islas Malvinas, now Montserrat (Montserrat overtakes Islas Malvinas)
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215 - was Montserrat, now Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
